$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters B, C, D, E, F, I, J, K, L, M map to these column indexes.
# (G and H are intentionally skipped: G stays 1 and H is empty, unchanged.)
$colIndexes = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13)

$newValues = @(
    @(1.02, 1.046410094356494, 1.052502357523927, 1.0535431997867, 1.063154658258426, 1.046370944666637, 1.051464639180589, 1.055250912285313, 1.056288878300932, 1.065874064711381),
    @(1.02, 1.047272777979757, 1.053174587009433, 1.054300472242773, 1.063969289289567, 1.046599402790631, 1.051975807316007, 1.055736361232914, 1.056859359196668, 1.066503652169135),
    @(1.02, 1.047831598638368, 1.053610051471106, 1.054791384495719, 1.064497371788747, 1.046746182611589, 1.052306500533726, 1.056050282978239, 1.057228738463213, 1.066911343259094),
    @(1.02, 1.048066670468117, 1.053793235560009, 1.054997979158089, 1.064719605997858, 1.046807637173305, 1.052445506687237, 1.056182207390799, 1.057384081523131, 1.067082808249699),
    @(1.02, 1.048106148438105, 1.053823999656795, 1.055032679874371, 1.06475693343283, 1.046817940882627, 1.052468845381368, 1.056204355210021, 1.057410167549512, 1.067111602132181),
    @(1.02, 1.047834739117329, 1.05361249873689, 1.054794144180609, 1.064500340398368, 1.046747004760378, 1.052308358010257, 1.056052045950098, 1.057230813946479, 1.06691363410311),
    @(1.02, 1.046701515740858, 1.052729438720105, 1.0537989347974, 1.063429766339915, 1.04644836963574, 1.051637403886412, 1.055415011856545, 1.056481624438279, 1.066086772452097),
    @(1.02, 1.044709349095275, 1.051177182344818, 1.052052270941792, 1.061550728715293, 1.045914147078352, 1.050454643086296, 1.054291032292906, 1.055163358215764, 1.06463214920438),
    @(1.02, 1.043384505847398, 1.05014501051104, 1.050892658676792, 1.060303160720004, 1.045552675682619, 1.049665907856401, 1.053540820388045, 1.054285877984772, 1.063664116507995),
    @(1.02, 1.042811629307165, 1.049698723224958, 1.050391702453034, 1.059764189429707, 1.045394902222928, 1.049324338948397, 1.053215775317032, 1.053906260887236, 1.063245374531638),
    @(1.02, 1.042598957648595, 1.049533051697715, 1.050205801679254, 1.059564178826632, 1.045336110543972, 1.049197460136541, 1.053095010631724, 1.053765306371851, 1.063089900200523),
    @(1.02, 1.042644570976324, 1.049568584273913, 1.050245670001556, 1.059607073271414, 1.04534873003468, 1.049224676273729, 1.053120916333812, 1.05379553920946, 1.063123247001521),
    @(1.02, 1.042794047340253, 1.049685026725868, 1.050376332222235, 1.059747652658051, 1.045390046306562, 1.049313851198888, 1.053205793443634, 1.053894608470234, 1.063232521638012),
    @(1.02, 1.042886160636683, 1.04975678399432, 1.050456860995384, 1.059834293139595, 1.045415477791653, 1.04936879420158, 1.053258085299052, 1.053955655268876, 1.063299858020016),
    @(1.02, 1.043422542544087, 1.050174642975127, 1.050925930165103, 1.060338956632884, 1.04556312023646, 1.049688575899783, 1.053562388493, 1.054311079158751, 1.063691916073265),
    @(1.02, 1.043759213203461, 1.050436930148016, 1.051220477691647, 1.060655850435328, 1.045655397286463, 1.049889156358778, 1.053753217677048, 1.054534118593506, 1.06393795799348),
    @(1.02, 1.043955663553817, 1.050589980406753, 1.051392394420272, 1.060840808149807, 1.045709099885683, 1.050006147470151, 1.053864505883019, 1.05466424620164, 1.064081510747569),
    @(1.02, 1.04402266089721, 1.05064217715117, 1.051451032546115, 1.060903894062277, 1.045727390528874, 1.050046037684732, 1.053902448983391, 1.054708621818111, 1.064130465375822),
    @(1.02, 1.043723083750032, 1.050408782729849, 1.051188863917025, 1.060621838414723, 1.04564550934654, 1.049867636386093, 1.053732745490892, 1.05451018520301, 1.063911555806288),
    @(1.02, 1.042750026945753, 1.049650734566517, 1.050337850562918, 1.059706250345568, 1.045377884866258, 1.049287591531174, 1.053180800025629, 1.053865433590175, 1.06320034117172),
    @(1.02, 1.042138923726729, 1.049174695826753, 1.049803807786457, 1.059131668597816, 1.045208533787225, 1.048922866539275, 1.052833605924925, 1.053460355259775, 1.062753549645546),
    @(1.02, 1.042462814574661, 1.049426997750403, 1.050086816366984, 1.059436161773843, 1.045298412570522, 1.049116216356041, 1.053017675196039, 1.05367506577964, 1.062990365986273),
    @(1.02, 1.043739408860065, 1.050421501143368, 1.051203148474722, 1.060637206615606, 1.045649977655065, 1.049877360348188, 1.053741996050176, 1.054520999571799, 1.063923485692763),
    @(1.02, 1.045223802978284, 1.05157801567309, 1.052502982012899, 1.062035610253581, 1.046053198193782, 1.050760460945029, 1.054581770971373, 1.055503927492446, 1.065007908582331)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $colIndexes.Count; $j++) {
        $col = $colIndexes[$j]
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
